$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Japanese translation column (column G) for rows 2-5.
$ws.Range("G2").Value = "縦隔の拡大"
$ws.Range("G3").Value = "狭い縦隔"
$ws.Range("G4").Value = "気胸"
$ws.Range("G5").Value = "縦隔リンパ節腫脹"

# Select the whole first row (header row), matching the new selection state.
$ws.Rows("1:1").Select() | Out-Null
